# Apply BOM updates for Thymio2 "2.1 output files" generation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: 10uF capacitor bank gains designator C78, quantity 9 -> 10
$ws.Range("A6").Value = 10
$ws.Range("C6").Value = "C2, C4, C31, C40, C46, C51, C52, C73, C76, C78"

# Row 24: inductor part swapped from MLP2520S3R3S to 2508056017Y2,
# designator list gains L5, quantity 4 -> 5
$ws.Range("A24").Value = 5
$ws.Range("B24").Value = "2508056017Y2"
$ws.Range("C24").Value = "L1, L2, L3, L4, L5"
$ws.Range("D24").Value = "SMD EMI Suppression Ferrite Bead WE-CBF, Z = 600 Ohm"
$ws.Range("E24").Value = "INDC2012X11L"
